$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "'1658736030850"
$ws.Range("B14").Value = "'6"
$ws.Range("C14").Value = "'0"
$ws.Range("D14").Value = "'8"
$ws.Range("E14").Value = "'4e"
$ws.Range("F14").Value = "'1"
$ws.Range("G14").Value = "'1"

$ws.Range("A15").Value = "'1658736242723"
$ws.Range("B15").Value = "'6"
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'8"
$ws.Range("E15").Value = "'14"
$ws.Range("F15").Value = "'4"
$ws.Range("G15").Value = "'12"

$ws.Range("A16").Value = "'1658760534237"
$ws.Range("B16").Value = "'9"
$ws.Range("C16").Value = "'0"
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'A"
$ws.Range("F16").Value = "'5"
$ws.Range("G16").Value = "'A"

$ws.Range("A17").Value = "'1658761359618"
$ws.Range("B17").Value = "'1"
$ws.Range("C17").Value = "'0"
$ws.Range("D17").Value = "'1, 2"
$ws.Range("E17").Value = "'4d"
$ws.Range("F17").Value = "'4"
$ws.Range("G17").Value = "'1"

$ws.Range("A18").Value = "'1658332363822"
$ws.Range("B18").Value = "'1"
$ws.Range("C18").Value = "'0"
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "'4d"
$ws.Range("F18").Value = "'4"
$ws.Range("G18").Value = "'1"

$ws.Range("A19").Value = "'1658332402648"
$ws.Range("B19").Value = "'1"
$ws.Range("C19").Value = "'0"
$ws.Range("D19").Value = "'0"
$ws.Range("E19").Value = "'4d"
$ws.Range("F19").Value = "'4"
$ws.Range("G19").Value = "'1"

$ws.Range("A20").Value = "'1658828890713"
$ws.Range("B20").Value = "'6"
$ws.Range("C20").Value = "'8b"
$ws.Range("D20").Value = "'0"
$ws.Range("E20").Value = "'4d"
$ws.Range("F20").Value = "'4"
$ws.Range("G20").Value = "'1"

$ws.Range("A21").Value = "'1658916484762"
$ws.Range("B21").Value = "'1"
$ws.Range("C21").Value = "'0"
$ws.Range("D21").Value = "'0"
$ws.Range("E21").Value = "'2"
$ws.Range("F21").Value = "'4"
$ws.Range("G21").Value = "'13"

$ws.Range("A22").Value = "'1658918012431"
$ws.Range("B22").Value = "'1"
$ws.Range("C22").Value = "'1"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'4c"
$ws.Range("F22").Value = "'4"
$ws.Range("G22").Value = "'1"

$ws.Range("A23").Value = "'1658920688970"
$ws.Range("B23").Value = "'4"
$ws.Range("C23").Value = "'0"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "'1"
$ws.Range("F23").Value = "'4"
$ws.Range("G23").Value = "'13"

$ws.Range("A24").Value = "'1659087654686"
$ws.Range("B24").Value = "'4"
$ws.Range("C24").Value = "'0"
$ws.Range("D24").Value = "'3, 9"
$ws.Range("E24").Value = "'5"
$ws.Range("F24").Value = "'4"
$ws.Range("G24").Value = "'13"

$ws.Range("A25").Value = "'1659344123097"
$ws.Range("B25").Value = "'1"
$ws.Range("C25").Value = "'0"
$ws.Range("D25").Value = "'0"
$ws.Range("E25").Value = "'4d"
$ws.Range("F25").Value = "'4"
$ws.Range("G25").Value = "'1"

$ws.Range("A26").Value = "'1659697068194"
$ws.Range("B26").Value = "'4"
$ws.Range("C26").Value = "'0"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "'4c"
$ws.Range("F26").Value = "'4"
$ws.Range("G26").Value = "'1"

$ws.Range("A27").Value = "'1660135979838"
$ws.Range("B27").Value = "'4"
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = "'4"
$ws.Range("E27").Value = "'2"
$ws.Range("F27").Value = "'4"
$ws.Range("G27").Value = "'13"

$ws.Range("A28").Value = "'1660292547743"
$ws.Range("B28").Value = "'9"
$ws.Range("C28").Value = "'11"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "'A"
$ws.Range("F28").Value = "'5"
$ws.Range("G28").Value = "'A"

$ws.Range("A29").Value = "'1660825903929"
$ws.Range("B29").Value = "'4"
$ws.Range("C29").Value = "'0"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "'9"
$ws.Range("F29").Value = "'4"
$ws.Range("G29").Value = "'ND"

$ws.Range("A30").Value = "'1661251238333"
$ws.Range("B30").Value = "'9"
$ws.Range("C30").Value = "'0"
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'A"
$ws.Range("F30").Value = "'5"
$ws.Range("G30").Value = "'A"

$ws.Range("A31").Value = "'1661853748366"
$ws.Range("B31").Value = "'1"
$ws.Range("C31").Value = "'11"
$ws.Range("D31").Value = "'0"
$ws.Range("E31").Value = "'4d"
$ws.Range("F31").Value = "'4"
$ws.Range("G31").Value = "'1"

$ws.Range("A32").Value = "'1663072313013"
$ws.Range("B32").Value = "'6"
$ws.Range("C32").Value = "'0"
$ws.Range("D32").Value = "'5"
$ws.Range("E32").Value = "'ND"
$ws.Range("G32").Value = "'ND"

$ws.Range("A33").Value = "'1663072962898"
$ws.Range("B33").Value = "'6"
$ws.Range("C33").Value = "'1"
$ws.Range("D33").Value = "'0"
$ws.Range("E33").Value = "'9"
$ws.Range("F33").Value = "'4"
$ws.Range("G33").Value = "'13"

$ws.Range("A34").Value = "'1664284436141"
$ws.Range("B34").Value = "'4"
$ws.Range("C34").Value = "'0"
$ws.Range("D34").Value = "'0"
$ws.Range("E34").Value = "'5"
$ws.Range("F34").Value = "'4"
$ws.Range("G34").Value = "'13"

$ws.Range("A35").Value = "'1664284677736"
$ws.Range("B35").Value = "'1"
$ws.Range("C35").Value = "'0"
$ws.Range("D35").Value = "'0"
$ws.Range("E35").Value = "'5"
$ws.Range("F35").Value = "'4"
$ws.Range("G35").Value = "'13"

$ws.Range("A36").Value = "'1664285824607"
$ws.Range("B36").Value = "'1"
$ws.Range("C36").Value = "'1"
$ws.Range("D36").Value = "'0"
$ws.Range("E36").Value = "'5"
$ws.Range("F36").Value = "'4"
$ws.Range("G36").Value = "'13"

$ws.Range("A37").Value = "'1664356543473"
$ws.Range("B37").Value = "'6"
$ws.Range("C37").Value = "'0"
$ws.Range("D37").Value = "'0"
$ws.Range("E37").Value = "'14"
$ws.Range("F37").Value = "'4"
$ws.Range("G37").Value = "'12"

$ws.Range("A38").Value = "'1664356559124"
$ws.Range("B38").Value = "'6"
$ws.Range("C38").Value = "'0"
$ws.Range("D38").Value = "'0"
$ws.Range("E38").Value = "'ND"
$ws.Range("G38").Value = "'ND"

$ws.Range("A39").Value = "'1664356571739"
$ws.Range("B39").Value = "'6"
$ws.Range("C39").Value = "'0"
$ws.Range("D39").Value = "'0"
$ws.Range("E39").Value = "'ND"
$ws.Range("G39").Value = "'ND"

$ws.Range("A40").Value = "'1664356649635"
$ws.Range("B40").Value = "'6"
$ws.Range("C40").Value = "'0"
$ws.Range("D40").Value = "'8"
$ws.Range("E40").Value = "'4e"
$ws.Range("F40").Value = "'1"
$ws.Range("G40").Value = "'1"

$ws.Range("A41").Value = "'1664358119728"
$ws.Range("B41").Value = "'6"
$ws.Range("C41").Value = "'6"
$ws.Range("D41").Value = "'0"
$ws.Range("E41").Value = "'9"
$ws.Range("F41").Value = "'4"
$ws.Range("G41").Value = "'ND"

$ws.Range("A42").Value = "'1664358454304"
$ws.Range("B42").Value = "'4"
$ws.Range("C42").Value = "'0"
$ws.Range("D42").Value = "'8"
$ws.Range("E42").Value = "'7"
$ws.Range("F42").Value = "'4"
$ws.Range("G42").Value = "'10"

$ws.Range("A43").Value = "'1666263595546"
$ws.Range("B43").Value = "'4"
$ws.Range("C43").Value = "'0"
$ws.Range("D43").Value = "'3"
$ws.Range("E43").Value = "'3a"
$ws.Range("F43").Value = "'4"
$ws.Range("G43").Value = "'13"

$ws.Range("A44").Value = "'1666255758576"
$ws.Range("B44").Value = "'1"
$ws.Range("C44").Value = "'0"
$ws.Range("D44").Value = "'0"
$ws.Range("E44").Value = "'4e"
$ws.Range("F44").Value = "'1"
$ws.Range("G44").Value = "'1"

$ws.Range("A45").Value = "'1666266085050"
$ws.Range("B45").Value = "'1"
$ws.Range("C45").Value = "'0"
$ws.Range("D45").Value = "'0"
$ws.Range("E45").Value = "'3a"
$ws.Range("F45").Value = "'4"
$ws.Range("G45").Value = "'13"

$ws.Range("A46").Value = "'1666267003464"
$ws.Range("B46").Value = "'1"
$ws.Range("C46").Value = "'11"
$ws.Range("D46").Value = "'0"
$ws.Range("E46").Value = "'4e"
$ws.Range("F46").Value = "'1"
$ws.Range("G46").Value = "'1"

$ws.Range("A47").Value = "'1668780056167"
$ws.Range("B47").Value = "'6"
$ws.Range("C47").Value = "'8b"
$ws.Range("D47").Value = "'0"
$ws.Range("E47").Value = "'14"
$ws.Range("F47").Value = "'4"
$ws.Range("G47").Value = "'12"

$ws.Range("A48").Value = "'1668781989003"
$ws.Range("B48").Value = "'6"
$ws.Range("C48").Value = "'0"
$ws.Range("D48").Value = "'0"
$ws.Range("E48").Value = "'ND"
$ws.Range("G48").Value = "'ND"

$ws.Range("A49").Value = "'1669281355185"
$ws.Range("B49").Value = "'4"
$ws.Range("C49").Value = "'0"
$ws.Range("D49").Value = "'3"
$ws.Range("E49").Value = "'3a"
$ws.Range("F49").Value = "'4"
$ws.Range("G49").Value = "'13"

$ws.Range("A50").Value = "'1669283492983"
$ws.Range("B50").Value = "'6"
$ws.Range("C50").Value = "'0"
$ws.Range("D50").Value = "'3"
$ws.Range("E50").Value = "'14"
$ws.Range("F50").Value = "'4"
$ws.Range("G50").Value = "'12"

$ws.Range("A51").Value = "'1669284735935"
$ws.Range("B51").Value = "'4"
$ws.Range("C51").Value = "'0"
$ws.Range("D51").Value = "'0"
$ws.Range("E51").Value = "'2"
$ws.Range("F51").Value = "'4"
$ws.Range("G51").Value = "'13"
